# Generate Report for Handback
#
# Populates the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns (F, G, H) for the two language
# report sheets (zh-cn, de-de) now that handback has completed, and
# flips the Status column from "Ready for handoff" to
# "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$srcUrl = "https://github.com/OpenLocalizationTest/oltest/blob/258de5c0980e6fc5e50ac0beff1e8421bd38454c/e2e/a.md"

$zhXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d2252a8dc0596d499e6eb277700c8f41a8c66bc5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$deXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$deXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/89b3e38536707199db6ac6dd1d6f4d0a6e5e70f0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

function Update-LanguageSheet($ws, [string]$xlfName, [string]$xlfUrl, [string]$handbackDateTime) {

    foreach ($row in 2, 3) {
        # F: Latest Target File -> a.md (hyperlink to source file)
        $fCell = $ws.Range("F$row")
        $fCell.Value = "a.md"
        $ws.Hyperlinks.Add($fCell, $srcUrl, [System.Type]::Missing, [System.Type]::Missing, "a.md") | Out-Null

        # G: Latest Handback File -> the handback xlf (hyperlink)
        $gCell = $ws.Range("G$row")
        $gCell.Value = $xlfName
        $ws.Hyperlinks.Add($gCell, $xlfUrl, [System.Type]::Missing, [System.Type]::Missing, $xlfName) | Out-Null

        # H: Latest Handback DateTime -> timestamp handback completed
        $ws.Range("H$row").Value = $handbackDateTime

        # C: Status -> now in sync, handed back
        $ws.Range("C$row").Value = "Handed back: in sync with en-US"
    }
}

$wsZh = $wb.Worksheets.Item("zh-cn")
Update-LanguageSheet $wsZh $zhXlfName $zhXlfUrl "2016-03-22 00:29:56"

$wsDe = $wb.Worksheets.Item("de-de")
Update-LanguageSheet $wsDe $deXlfName $deXlfUrl "2016-03-22 00:30:06"

# The Overview sheet rolls up the per-language Status for each file via
# the same shared "Status" text, so it also flips to the handed-back
# wording for both zh-cn (col B) and de-de (col C).
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($row in 2, 3) {
    $wsOverview.Range("B$row").Value = "Handed back: in sync with en-US"
    $wsOverview.Range("C$row").Value = "Handed back: in sync with en-US"
}
